$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update an existing data point: Bayou Segnette Closure (row 8) gage now reading 2.5 ---
$ws.Range("C8").Value = 2.5

# --- Append three new gate/station rows (31-33), copying the style of the
#     preceding "extra stations" block (rows 23-30: style index 3 on col A,
#     style index 2 on col C) so no new cellXfs are introduced ---
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31:A33").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C31:C33").PasteSpecial(-4122) | Out-Null

$ws.Range("A31").Value = "Venice"
$ws.Range("B31").Value = "MS River at Venice"
$ws.Range("C31").Value = 0

$ws.Range("A32").Value = "BayouSale"
$ws.Range("B32").Value = "GIWW at Bayou Sale Ridge"
$ws.Range("C32").Value = 0

$ws.Range("A33").Value = "BayouBoeuf"
$ws.Range("B33").Value = "Bayou Boeuf at Railroad Bridge"
$ws.Range("C33").Value = 0

# --- Update the view state: scroll back to the top-left (drop the saved
#     topLeftCell="A2") and move the active selection to the new first
#     blank row below the table ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A34").Select() | Out-Null
